$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SYF")

$ws.Range("D12").Value = 1.0491
$ws.Range("E12").Value = 0.9757
$ws.Range("F12").Value = 0.9272
$ws.Range("G12").Value = 0.8823

$ws.Range("D14").Value = 0.1224
$ws.Range("E14").Value = 0.166
$ws.Range("F14").Value = 0.2096
$ws.Range("G14").Value = 0.2511

$ws.Range("D15").Value = 0.0929
$ws.Range("E15").Value = 0.1273
$ws.Range("F15").Value = 0.1605
$ws.Range("G15").Value = 0.1925

$ws.Range("D16").Value = 0.534
$ws.Range("E16").Value = 0.5249
$ws.Range("F16").Value = 0.4536
$ws.Range("G16").Value = 0.4619

$ws.Range("D23").Value = 0.534
$ws.Range("E23").Value = 0.5249
$ws.Range("F23").Value = 0.4536
$ws.Range("G23").Value = 0.4619
